$d = $word.ActiveDocument

# Locate the end of the text "¿Qué servicios ofreceré?" (this point sits right before
# the paragraph mark that currently also carries the "_GoBack" bookmark).
$find1 = $d.Content
$found1 = $find1.Find.Execute("¿Qué servicios ofreceré?", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)

# Insert the new "ARQUITECTURA POR CAPAS" paragraph (plus its trailing paragraph mark)
# right before that point, so the bookmark stays attached to the paragraph mark that
# follows - i.e. it ends up owning its own (now empty) paragraph.
$insertPoint1 = $d.Range($find1.End, $find1.End)
$insertPoint1.InsertBefore("ARQUITECTURA POR CAPAS`r")

# Split "¿Qué servicios ofreceré?" away from the text that was just inserted, giving us a
# dedicated paragraph for the original sentence.
$find2 = $d.Content
$found2 = $find2.Find.Execute("ofreceré?", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
$insertPoint2 = $d.Range($find2.End, $find2.End)
$insertPoint2.InsertBefore("`r")

# Add the blank paragraph that separates "¿Qué servicios ofreceré?" from
# "ARQUITECTURA POR CAPAS".
$find3 = $d.Content
$found3 = $find3.Find.Execute("ARQUITECTURA POR CAPAS", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
$insertPoint3 = $d.Range($find3.Start, $find3.Start)
$insertPoint3.InsertBefore("`r")

# Center the new heading paragraph.
$find4 = $d.Content
$found4 = $find4.Find.Execute("ARQUITECTURA POR CAPAS", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
$find4.Paragraphs.Item(1).Format.Alignment = 1
